# Update "想去人数" (want-to-go count) figures and one event title across all
# four sheets of the workbook, per the upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 4624
$ws1.Range("F8").Value  = 721
$ws1.Range("F12").Value = 1126
$ws1.Range("F14").Value = 804
$ws1.Range("F15").Value = 759
$ws1.Range("F16").Value = 546
$ws1.Range("F17").Value = 509
$ws1.Range("F19").Value = 157
$ws1.Range("F22").Value = 386
$ws1.Range("F23").Value = 2499
$ws1.Range("F25").Value = 1531
$ws1.Range("F26").Value = 485
$ws1.Range("F29").Value = 4217

# --- Sheet 2: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value  = 359
$ws2.Range("F15").Value = 49
$ws2.Range("F17").Value = 280
$ws2.Range("F18").Value = 279
$ws2.Range("C19").Value = "上海·【早鸟5折】红楼梦·梁祝·探清水河 《国潮》跨界音乐会"
$ws2.Range("F20").Value = 139

# --- Sheet 3: 本地生活 (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 99
$ws3.Range("F4").Value = 1325
$ws3.Range("F5").Value = 1721
$ws3.Range("F6").Value = 1069
$ws3.Range("F7").Value = 229

# --- Sheet 4: 全部类型 (All Types - combined view) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1325
$ws4.Range("F3").Value  = 1721
$ws4.Range("F4").Value  = 1069
$ws4.Range("F5").Value  = 229
$ws4.Range("F8").Value  = 4624
$ws4.Range("F12").Value = 721
$ws4.Range("F13").Value = 359
$ws4.Range("F17").Value = 1126
$ws4.Range("F21").Value = 804
$ws4.Range("F22").Value = 759
$ws4.Range("F23").Value = 546
$ws4.Range("F24").Value = 509
$ws4.Range("F26").Value = 157
$ws4.Range("F27").Value = 49
$ws4.Range("F28").Value = 280
$ws4.Range("F29").Value = 279
$ws4.Range("C30").Value = "上海·【早鸟5折】红楼梦·梁祝·探清水河 《国潮》跨界音乐会"
$ws4.Range("C31").Value = "上海·【早鸟5折】红楼梦·梁祝·探清水河 《国潮》跨界音乐会"
$ws4.Range("F34").Value = 386
$ws4.Range("F35").Value = 139
$ws4.Range("F37").Value = 2499
$ws4.Range("F43").Value = 1532
$ws4.Range("F44").Value = 485
$ws4.Range("F48").Value = 4217
